$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.262.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.677.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.05%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.45%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.36%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.105"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.40%  "

# Row 10
$ws.Range("E10").Value = "  -0.53%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.367"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.11%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.64%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.151.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.94%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.131.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.05%  "

# Row 16
$ws.Range("E16").Value = "  -4.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.679.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.99%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.49%  "

# Row 22
$ws.Range("E22").Value = "  -0.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.505"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.59%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0849"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.45%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.81%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.60%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.37%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.83%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.85%  "

# Row 36
$ws.Range("E36").Value = "  -5.58%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "338.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.939"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.72%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.04%  "

# Row 41
$ws.Range("E41").Value = "  -2.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.11%  "

# Row 43
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.22%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.91%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.617"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.83%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0558"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.22%  "

# Row 47
$ws.Range("E47").Value = "  +0.11%  "

# Row 48
$ws.Range("E48").Value = "  +0.00%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.29%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0967"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.24%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.30%  "
